$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B30" = "e1abb237b7653418ea1339950d90d2be"
    "B34" = "fb314a97514b3e24c4a895f6dab99cd3"
    "B57" = "bb2240a3f8966e232606fd3fabe29348"
    "B115" = "78fb34603fc974bb8815be6ff28d67f3"
    "B117" = "c44933a8687ca715bd1e53da6d63de28"
    "B154" = "0164192226833e8b2508d9634b0ba903"
    "B161" = "1e5c3f3bf56fea72588394470e1cc359"
    "B191" = "dfcc9d17c7339038264c5f0bce129efc"
    "B246" = "98d307afdb21d3f40b972bb11077df14"
    "B281" = "825ad91b8f396f7b2ff56467cabb4d7a"
    "B299" = "27b7354351f85b3ec9741b3dc249118a"
    "B387" = "be634ffb672776a9bf2fb361968e241c"
    "B410" = "37880be478263c4ee5741708a79ffd3d"
    "B415" = "5bc66926ec0893680b606c0d50c3c2cd"
    "B419" = "0841f66eec1f7caf51680bed6f5054c6"
    "B424" = "8fef2a16bee470f5ea90550ee9ece9e0"
    "B480" = "23abdb30cf3035023e1aef078106e96e"
    "B504" = "2c7c22ed1ce7767e03ff2638310fc76b"
    "B524" = "5756fa8a1ed6e839d07823f1691edd57"
    "B552" = "26800b7f7072242058a59dc30201fb2a"
    "B601" = "122e50541cdb47f369c40eb3484e3e6d"
    "B618" = "cdeec3a4e361cc7e3e633c7a2be1280d"
    "B649" = "c57bf5965e73c6bcb8e711980866402a"
    "B655" = "3379e70f93178a55f709d366d220e3ba"
    "B666" = "2da3844e6a18aaa5b5a36d9d8baf282a"
    "B704" = "2a3504935d9b2f689225d9b3cdd46f48"
    "B712" = "112d140bc3db4e5bca98e0b3941c811b"
    "B729" = "ab4a1ff81ca4806a30c66cda47c56e6d"
    "B740" = "67941c3e10fa100f277461489faf20df"
    "B742" = "47db683e5277c6d835ffe16eae4a2371"
    "B756" = "f0eff0ceb9a47a301a34844e3837a4a3"
    "B785" = "36609a79c07265c4ab1161460dbbdf6c"
    "B786" = "20a6ca4e79ebeafd2fa55b0026034986"
    "B802" = "5b60f5dc7d5d4b86bcb9fb200e716a38"
    "B811" = "5f1e48ea2ee37ac4a0cd6534daf28e1d"
    "B816" = "831b12f239db1883cfb6a62cd480eabe"
    "B830" = "e201a0c7853b150c6ba269bc741a519b"
    "B839" = "f45609a4e0bac1efd6962c8af389dc06"
    "B846" = "19d9b73b5cc05a07fd97c34a5ce6f55a"
    "B848" = "ae7efecaf8736ca69f95c36d2f77d0d1"
    "B874" = "c9c849f03081bb7a17b5eba5feebb7ea"
    "B911" = "00bbac0f63cee336177391fe8fd966bc"
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
